# Auto-generated Excel COM-interop script to apply cryptos.xlsx update
# Commit message: Updated cryptos list on Fri Jun  2 10:48:27 UTC 2023 with GitHub Actions
#
# Every changed cell in the source workbook is stored as inline/shared TEXT
# (t="inlineStr"), even when the text looks like a number (e.g. "1.000",
# "0.5147", prices with thousand-dot separators, etc). When such a string is
# written through Range.Value, Excel auto-detects it as a number and converts
# the cell to a numeric cell, which must be avoided. To keep these cells as
# plain text (matching the original file) we temporarily force a Text number
# format before assigning the value, then restore the "Normal" style so the
# cell keeps its original (default) appearance/format afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.081.14"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "1.890.45"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5147"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.94%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3768"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07214"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("E10").Value = "  +2.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9042"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07651"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.34%  "
$ws.Range("D13").Value = "1.883.49"
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.32%  "
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008475"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9997"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "27.103.46"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.068"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("D22").Value = "2.126.54"
$ws.Range("E22").Value = "  +2.05%  "
$ws.Range("E23").Value = "  +2.00%  "
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.37%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.247"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.19%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.781"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.948"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.829"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09179"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05084"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.235"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7828"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.002"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.287"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("E38").Value = "  +4.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.074"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.626"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "117.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.70%  "
$ws.Range("E45").Value = "  +2.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4800"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9997"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.595"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.83%  "
$ws.Range("E50").Value = "  +1.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.62%  "
